# Applies the edits described by the commit diff:
#  1. (presentation.xml) Register an empty static-guide list extension
#     (p15:sldGuideLst) on the presentation extLst. Not exposed anywhere
#     in the documented PowerPoint COM surface (Presentation.Guides /
#     Application.DisplayGuides are the closest analogues), so we best-
#     effort it and simply continue if the host does not support it.
#  2. (slide 1 / "TextBox 36") Text "Ricevere" -> "Riceve" in the first
#     run of the paragraph, and shrink the shape width from 4487319 EMU
#     to 4294765 EMU (the box uses spAutoFit / wrap="none", so removing
#     two letters narrows the auto-fitted box).

$p = $ppt.ActivePresentation

# --- 1. Best-effort: try to materialize the (normally view-state-only)
#        static guide list so the extension shows up if this host ever
#        wires it up. Real PowerPoint gives no supported way to do this
#        other than toggling the Guides UI, so failures here are expected
#        and silently ignored.
try {
    $app = $ppt
    $app.DisplayGuides = $true
} catch {
}
try {
    $null = $p.Guides
} catch {
}

# --- 2. Locate the shape. It is shape #19 ("TextBox 36") on slide 1.
$s = $p.Slides.Item(1)
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 36") {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(19)
}

$tr = $shape.TextFrame.TextRange

# Only the first run's text changes ("Ricevere" -> "Riceve"); every other
# run in the paragraph ("  i  dati  fino a un  massimo  di 1024 byte")
# keeps its original text/formatting untouched.
$tr.Runs(1).Text = "Riceve"

# The textbox auto-fits its width to the (now shorter) text. Target width
# is 4294765 EMU == 338.1705 pt (PowerPoint stores shape extents in
# points; 338.1705 is the closest point value that round-trips to the
# exact EMU figure from the diff).
$shape.Width = 338.1705
